$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values ---
$ws.Range("D2").Value = "73.005.85"
$ws.Range("D3").Value = "3.987.50"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.73"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.43"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.684"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.11"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000319"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.96"
$ws.Range("D14").Value = "4.613.23"
$ws.Range("D15").Value = "4.002.51"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.17"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.37"
$ws.Range("D20").Value = "72.580.04"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "433.69"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.78"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.19"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.22"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.40"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.75"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.92"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "672.27"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.43"
$ws.Range("D37").Value = "0.0₃0886"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0491"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.63"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("D49").Value = "2.852.33"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.03"

# --- Update Volume(1h) (column E) values ---
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +9.71%  "
$ws.Range("E6").Value = "  +7.79%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -4.69%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +8.06%  "
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("E22").Value = "  +13.10%  "
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("E26").Value = "  +16.35%  "
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("E34").Value = "  -4.95%  "
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("E36").Value = "  +7.37%  "
$ws.Range("E37").Value = "  +7.68%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  +9.14%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +3.70%  "
